$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the contiguous block of rows 2-11 (completed demand entries),
# shifting the remaining rows up.
$ws.Range("A2:A11").EntireRow.Delete()

# Update the active selection to match the new data range
$ws.Range("A2:A28").Select()
